$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.440985666666667
$ws.Range("H2").Value = 4.322957000000001
$ws.Range("I2").Value = 0.1098365531732288
$ws.Range("J2").Value = 0.1230162332390494
$ws.Range("M2").Value = 6.712486666666666
$ws.Range("N2").Value = 20.13746
$ws.Range("O2").Value = 0.6330487633990675
$ws.Range("P2").Value = 0.6414503882251803
$ws.Range("Q2").Value = 9.672597074357778
$ws.Range("R2").Value = 87.05337366921999
$ws.Range("S2").Value = 0.06953189416232841
$ws.Range("T2").Value = 0.07890881056918754
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.440985666666667
$ws.Range("H3").Value = 4.322957000000001
$ws.Range("I3").Value = 0.1098365531732288
$ws.Range("J3").Value = 0.1230162332390494
$ws.Range("O3").Value = 0.290741083484562
$ws.Range("P3").Value = 0.2945997080427384
$ws.Range("Q3").Value = 4.442345544454668
$ws.Range("R3").Value = 39.98110990009201
$ws.Range("S3").Value = 0.03193399847579424
$ws.Range("T3").Value = 0.03624054639674135
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.440985666666667
$ws.Range("H4").Value = 4.322957000000001
$ws.Range("I4").Value = 0.1098365531732288
$ws.Range("J4").Value = 0.1230162332390494
$ws.Range("M4").Value = 0.2495096666666667
$ws.Range("N4").Value = 0.748529
$ws.Range("O4").Value = 0.02353103905946135
$ws.Range("P4").Value = 0.02384333563656022
$ws.Range("Q4").Value = 0.3595398533614445
$ws.Range("R4").Value = 3.235858680253
$ws.Range("S4").Value = 0.002584568222875851
$ws.Range("T4").Value = 0.00293311733786403
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.440985666666667
$ws.Range("H5").Value = 4.322957000000001
$ws.Range("I5").Value = 0.1098365531732288
$ws.Range("J5").Value = 0.1230162332390494
$ws.Range("M5").Value = 0.4166465
$ws.Range("N5").Value = 0.8332930000000001
$ws.Range("O5").Value = 0.03929356804674715
$ws.Range("P5").Value = 0.02654337331298611
$ws.Range("Q5").Value = 0.6003816345668335
$ws.Range("R5").Value = 3.602289807401001
$ws.Range("S5").Value = 0.004315870076132427
$ws.Range("T5").Value = 0.003265265802421458
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.440985666666667
$ws.Range("H6").Value = 4.322957000000001
$ws.Range("I6").Value = 0.1098365531732288
$ws.Range("J6").Value = 0.1230162332390494
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1419326666666667
$ws.Range("N6").Value = 0.425798
$ws.Range("O6").Value = 0.01338554601016197
$ws.Range("P6").Value = 0.01356319478253491
$ws.Range("Q6").Value = 0.2045229382984445
$ws.Range("R6").Value = 1.840706444686
$ws.Range("S6").Value = 0.001470222236097855
$ws.Range("T6").Value = 0.001668493132834971
$ws.Range("I7").Value = 0.5687502547919595
$ws.Range("J7").Value = 0.6369966279614609
$ws.Range("M7").Value = 6.712486666666666
$ws.Range("N7").Value = 20.13746
$ws.Range("O7").Value = 0.6330487633990675
$ws.Range("P7").Value = 0.6414503882251803
$ws.Range("Q7").Value = 50.08616796144888
$ws.Range("R7").Value = 450.7755116530399
$ws.Range("S7").Value = 0.3600466454789545
$ws.Range("T7").Value = 0.4086017343040099
$ws.Range("I8").Value = 0.5687502547919595
$ws.Range("J8").Value = 0.6369966279614609
$ws.Range("O8").Value = 0.290741083484562
$ws.Range("P8").Value = 0.2945997080427384
$ws.Range("S8").Value = 0.165359065310335
$ws.Range("T8").Value = 0.1876590206216552
$ws.Range("I9").Value = 0.5687502547919595
$ws.Range("J9").Value = 0.6369966279614609
$ws.Range("M9").Value = 0.2495096666666667
$ws.Range("N9").Value = 0.748529
$ws.Range("O9").Value = 0.02353103905946135
$ws.Range("P9").Value = 0.02384333563656022
$ws.Range("Q9").Value = 1.861751641866222
$ws.Range("R9").Value = 16.755764776796
$ws.Range("S9").Value = 0.01338328446058819
$ws.Range("T9").Value = 0.0151881243998422
$ws.Range("I10").Value = 0.5687502547919595
$ws.Range("J10").Value = 0.6369966279614609
$ws.Range("M10").Value = 0.4166465
$ws.Range("N10").Value = 0.8332930000000001
$ws.Range("O10").Value = 0.03929356804674715
$ws.Range("P10").Value = 0.02654337331298611
$ws.Range("Q10").Value = 3.108866745788667
$ws.Range("R10").Value = 18.653200474732
$ws.Range("S10").Value = 0.02234822683827264
$ws.Range("T10").Value = 0.01690803929509439
$ws.Range("I11").Value = 0.5687502547919595
$ws.Range("J11").Value = 0.6369966279614609
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1419326666666667
$ws.Range("N11").Value = 0.425798
$ws.Range("O11").Value = 0.01338554601016197
$ws.Range("P11").Value = 0.01356319478253491
$ws.Range("Q11").Value = 1.059050652150222
$ws.Range("R11").Value = 9.531455869352
$ws.Range("S11").Value = 0.007613032703809116
$ws.Range("T11").Value = 0.008639709340859215
$ws.Range("G12").Value = 4.2167365
$ws.Range("H12").Value = 8.433472999999999
$ws.Range("I12").Value = 0.3214131920348118
$ws.Range("J12").Value = 0.2399871387994896
$ws.Range("M12").Value = 6.712486666666666
$ws.Range("N12").Value = 20.13746
$ws.Range("O12").Value = 0.6330487633990675
$ws.Range("P12").Value = 0.6414503882251803
$ws.Range("Q12").Value = 28.30478753309666
$ws.Range("R12").Value = 169.82872519858
$ws.Range("S12").Value = 0.2034702237577846
$ws.Range("T12").Value = 0.1539398433519828
$ws.Range("G13").Value = 4.2167365
$ws.Range("H13").Value = 8.433472999999999
$ws.Range("I13").Value = 0.3214131920348118
$ws.Range("J13").Value = 0.2399871387994896
$ws.Range("O13").Value = 0.290741083484562
$ws.Range("P13").Value = 0.2945997080427384
$ws.Range("Q13").Value = 12.999574552498
$ws.Range("R13").Value = 77.99744731498799
$ws.Range("S13").Value = 0.09344801969843275
$ws.Range("T13").Value = 0.07070014102434177
$ws.Range("G14").Value = 4.2167365
$ws.Range("H14").Value = 8.433472999999999
$ws.Range("I14").Value = 0.3214131920348118
$ws.Range("J14").Value = 0.2399871387994896
$ws.Range("M14").Value = 0.2495096666666667
$ws.Range("N14").Value = 0.748529
$ws.Range("O14").Value = 0.02353103905946135
$ws.Range("P14").Value = 0.02384333563656022
$ws.Range("Q14").Value = 1.052116518536167
$ws.Range("R14").Value = 6.312699111216999
$ws.Range("S14").Value = 0.007563186375997308
$ws.Range("T14").Value = 0.005722093898853996
$ws.Range("G15").Value = 4.2167365
$ws.Range("H15").Value = 8.433472999999999
$ws.Range("I15").Value = 0.3214131920348118
$ws.Range("J15").Value = 0.2399871387994896
$ws.Range("M15").Value = 0.4166465
$ws.Range("N15").Value = 0.8332930000000001
$ws.Range("O15").Value = 0.03929356804674715
$ws.Range("P15").Value = 0.02654337331298611
$ws.Range("Q15").Value = 1.75688850414725
$ws.Range("R15").Value = 7.027554016589
$ws.Range("S15").Value = 0.01262947113234208
$ws.Range("T15").Value = 0.006370068215470266
$ws.Range("G16").Value = 4.2167365
$ws.Range("H16").Value = 8.433472999999999
$ws.Range("I16").Value = 0.3214131920348118
$ws.Range("J16").Value = 0.2399871387994896
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1419326666666667
$ws.Range("N16").Value = 0.425798
$ws.Range("O16").Value = 0.01338554601016197
$ws.Range("P16").Value = 0.01356319478253491
$ws.Range("Q16").Value = 0.5984926560756667
$ws.Range("R16").Value = 3.590955936454
$ws.Range("S16").Value = 0.004302291070254997
$ws.Range("T16").Value = 0.003254992308840718
